$d = $word.ActiveDocument

# 1. Fill in the Github project link placeholder.
$d.Content.Find.Execute(
    "Lien Github du projet : ……………",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Lien Github du projet : https://github.com/DroxKiwi/evaluation-2.git",
    2
)

# 2. Fill in the three previously empty paragraphs of the "Décrivez les
#    tâches..." answer box with the author's write-up.
$paras = $d.Paragraphs

$p1 = $paras.Item(28)
$p1.Range.Text = "C’est un site vitrine, qui a pour but de mettre en avant un peu de JavaScript comme demandé pour l’évaluation. "

$p2 = $paras.Item(29)
$p2.Range.Text = "J’ai d’abord repris le contenu de mon évaluation 1 et ajouté une page js.html au projet. "

$p3 = $paras.Item(30)
$p3.Range.Text = "La page contient une « div » qui est ajustable en largeur. Lorsque la div est agrandi, la page html est modifié comme si un filtre était appliqué."
